# ==========================================================================
# [ADDITIONAL SCRAPING] add a "Player Info" sheet and an "ODI Batting Extra"
# sheet, and replace the full match-card-URL column on the existing
# "ODI Batting" / "ODI Bowling" sheets with a short numeric MATCH_CODE.
# ==========================================================================

$wb = $excel.ActiveWorkbook

# --------------------------------------------------------------------
# Helper: write a row of values into a worksheet starting at (row, col).
# Any value that is a [string] gets a leading apostrophe so Excel keeps
# it as literal text instead of silently re-typing it as a number /
# percentage (mirrors how a human/VBA author would force text entry).
# --------------------------------------------------------------------
function Set-TextRow {
    param($ws, [int]$row, [int]$startCol, [object[]]$values)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $val = $values[$i]
        if ($null -eq $val) { continue }
        $ws.Cells.Item($row, $startCol + $i).Value = "'" + $val
    }
}

# ==========================================================================
# 1) Existing sheet "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE
# ==========================================================================
$wsBatting = $wb.Worksheets.Item("ODI Batting")

$wsBatting.Range("D1").Value = "MATCH_CODE"

$battingCodes = @(
    "3656","3665","3671","3685","3686","3687","3695","3697","3700","3728",
    "3732","3736","3740","3742","3766","3771","3776","3783","3790","3794",
    "3813","3816","3818","3826","3827","3828","3885","3887","3891","3892",
    "3894","3940","3942","3945","3947","3950"
)
for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $r = $i + 2
    $wsBatting.Cells.Item($r, 4).Value = "'" + $battingCodes[$i]
}

# ==========================================================================
# 2) Existing sheet "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE
# ==========================================================================
$wsBowling = $wb.Worksheets.Item("ODI Bowling")

$wsBowling.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @("3671","3771","3827")
for ($i = 0; $i -lt $bowlingCodes.Length; $i++) {
    $r = $i + 2
    $wsBowling.Cells.Item($r, 2).Value = "'" + $bowlingCodes[$i]
}

# ==========================================================================
# 3) New sheet "Player Info", inserted before "ODI Batting"
# ==========================================================================
$wsPlayer = $wb.Worksheets.Add($wsBatting)
$wsPlayer.Name = "Player Info"

$playerHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 0; $c -lt $playerHeaders.Length; $c++) {
    $cell = $wsPlayer.Cells.Item(1, $c + 1)
    $cell.Value = $playerHeaders[$c]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

$wsPlayer.Cells.Item(2, 1).Value = "'4321"
$wsPlayer.Cells.Item(2, 2).Value = "Rilee Roscoe Rossouw"
$wsPlayer.Cells.Item(2, 3).Value = "Left Handed"
$wsPlayer.Cells.Item(2, 4).Value = "Right Arm Off Break"

# ==========================================================================
# 4) New sheet "ODI Batting Extra", inserted after "ODI Bowling"
#    (re-fetch the "ODI Bowling" reference - the sheet collection shifted
#    when "Player Info" was inserted above, so the old handle is stale)
# ==========================================================================
$wsBowlingFresh = $wb.Worksheets.Item("ODI Bowling")
$wsExtra = $wb.Worksheets.Add($null, $wsBowlingFresh)
$wsExtra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 0; $c -lt $extraHeaders.Length; $c++) {
    $cell = $wsExtra.Cells.Item(1, $c + 1)
    $cell.Value = $extraHeaders[$c]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("3776", 4,    "1",  "0", "2.97%",  "NO"),
    @("3783", 3,    "8",  "0", "12.61%", "NO"),
    @("3790", 4,    $null,$null,$null,   "NO"),
    @("3794", 4,    "2",  "1", "13.88%", "NO"),
    @("3813", 4,    "5",  "1", "27.44%", "NO"),
    @("3816", 4,    "0",  "0", "2.47%",  "NO"),
    @("3818", 4,    "2",  "0", "10.12%", "NO"),
    @("3826", 3,    "6",  "3", "29.28%", "NO"),
    @("3827", 3,    "5",  "0", "19.12%", "NO"),
    @("3828", 3,    "0",  "0", "2.12%",  "NO"),
    @("3885", 6,    "3",  "0", "7.60%",  "NO"),
    @("3887", 6,    "2",  "0", "4.20%",  "NO"),
    @("3891", 4,    "1",  "0", "1.69%",  "NO"),
    @("3892", 3,    "2",  "0", "32.45%", "NO"),
    @("3894", 3,    "1",  "0", "3.70%",  "NO"),
    @("3940", $null,$null,$null,$null,   "NO"),
    @("3942", 2,    "10", "0", "20.78%", "NO"),
    @("3945", 4,    "1",  "0", "4.84%",  "NO"),
    @("3947", $null,$null,$null,$null,   "NO"),
    @("3950", $null,$null,$null,$null,   "NO")
)

for ($i = 0; $i -lt $extraRows.Length; $i++) {
    $r = $i + 2
    $row = $extraRows[$i]

    # A: MATCH_CODE - text
    $wsExtra.Cells.Item($r, 1).Value = "'" + $row[0]

    # B: BATTING_POSITION - numeric
    if ($null -ne $row[1]) {
        $wsExtra.Cells.Item($r, 2).Value = $row[1]
    }

    # C, D: NUM_4 / NUM_6 - text
    if ($null -ne $row[2]) { $wsExtra.Cells.Item($r, 3).Value = "'" + $row[2] }
    if ($null -ne $row[3]) { $wsExtra.Cells.Item($r, 4).Value = "'" + $row[3] }

    # E: PERCENT_RUNS_OF_TOTAL - text
    if ($null -ne $row[4]) { $wsExtra.Cells.Item($r, 5).Value = "'" + $row[4] }

    # F: MAN_OF_MATCH - text (safe, not numeric-looking)
    $wsExtra.Cells.Item($r, 6).Value = $row[5]
}
